$wb = $excel.ActiveWorkbook

# Mapping of row number -> new "想去人数" (F column) value.
# Same updates apply to both the "展览" and "全部类型" worksheets,
# which contain identical data.
$updates = @{
    2  = 3146
    4  = 166
    6  = 1769
    8  = 105
    11 = 1464
    14 = 362
    15 = 87
    17 = 84
    20 = 131
    23 = 3447
    24 = 416
    25 = 300
    26 = 473
    27 = 79
    28 = 24
    30 = 1197
    31 = 127
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
